$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '27.191.37'
$ws.Range('E2').Value = '  -2.38%  '
$ws.Range('D3').Value = '1.720.87'
$ws.Range('E3').Value = '  -2.56%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextCell 'D5' '312.64'
$ws.Range('E5').Value = '  -4.46%  '
Set-TextCell 'D6' '1.002'
$ws.Range('E6').Value = '  +0.02%  '
Set-TextCell 'D7' '0.4609'
$ws.Range('E7').Value = '  +3.19%  '
Set-TextCell 'D8' '0.3444'
$ws.Range('E8').Value = '  -2.69%  '
Set-TextCell 'D9' '42.48'
$ws.Range('E9').Value = '  +1.08%  '
Set-TextCell 'D10' '0.07268'
$ws.Range('E10').Value = '  -2.20%  '
Set-TextCell 'D11' '1.042'
$ws.Range('E11').Value = '  -4.61%  '
Set-TextCell 'D13' '19.80'
$ws.Range('E13').Value = '  -4.79%  '
Set-TextCell 'D14' '5.844'
$ws.Range('E14').Value = '  -2.85%  '
$ws.Range('D15').Value = '1.724.31'
$ws.Range('E15').Value = '  -2.59%  '
Set-TextCell 'D16' '6.872'
$ws.Range('E16').Value = '  -4.42%  '
Set-TextCell 'D17' '89.52'
$ws.Range('E17').Value = '  -3.57%  '
$ws.Range('E18').Value = '  -1.69%  '
Set-TextCell 'D19' '0.06330'
$ws.Range('E19').Value = '  -1.52%  '
Set-TextCell 'D20' '1.001'
$ws.Range('E20').Value = '  +0.01%  '
Set-TextCell 'D21' '16.49'
$ws.Range('E21').Value = '  -3.72%  '
Set-TextCell 'D22' '5.623'
$ws.Range('E22').Value = '  -2.65%  '
$ws.Range('D23').Value = '27.232.84'
$ws.Range('E23').Value = '  -2.39%  '
Set-TextCell 'D24' '10.87'
$ws.Range('E24').Value = '  -3.53%  '
Set-TextCell 'D25' '2.133'
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('E26').Value = '  -4.85%  '
Set-TextCell 'D27' '19.30'
$ws.Range('E27').Value = '  -4.20%  '
$ws.Range('D28').Value = '1.921.25'
$ws.Range('E28').Value = '  -2.67%  '
Set-TextCell 'D29' '2.154'
$ws.Range('E29').Value = '  -0.51%  '
Set-TextCell 'D30' '119.10'
$ws.Range('E30').Value = '  -4.59%  '
Set-TextCell 'D31' '1.028'
$ws.Range('E31').Value = '  -6.16%  '
Set-TextCell 'D32' '0.09091'
$ws.Range('E32').Value = '  -0.51%  '
Set-TextCell 'D33' '3.593'
$ws.Range('E33').Value = '  -1.50%  '
Set-TextCell 'D34' '5.345'
$ws.Range('E34').Value = '  -3.70%  '
Set-TextCell 'D35' '0.02207'
$ws.Range('E35').Value = '  -3.60%  '
Set-TextCell 'D36' '0.05860'
$ws.Range('E36').Value = '  -3.84%  '
Set-TextCell 'D37' '11.08'
$ws.Range('E37').Value = '  -6.22%  '
Set-TextCell 'D38' '0.1998'
$ws.Range('E38').Value = '  -4.29%  '
$ws.Range('B39').Value = 'WEMIXTOKEN'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 'D39' '1.420'
$ws.Range('E39').Value = '  +2.36%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D40' '4.729'
$ws.Range('E40').Value = '  -4.75%  '
Set-TextCell 'D41' '0.5940'
$ws.Range('E41').Value = '  -5.59%  '
$ws.Range('E42').Value = '  -4.06%  '
Set-TextCell 'D43' '7.465'
$ws.Range('E43').Value = '  -5.60%  '
Set-TextCell 'D44' '12.79'
$ws.Range('E44').Value = '  -3.28%  '
Set-TextCell 'D45' '3.594'
$ws.Range('E45').Value = '  -3.70%  '
Set-TextCell 'D46' '0.5637'
$ws.Range('E46').Value = '  -3.46%  '
Set-TextCell 'D47' '119.83'
$ws.Range('E47').Value = '  -1.70%  '
Set-TextCell 'D48' '1.868'
$ws.Range('E48').Value = '  -3.94%  '
Set-TextCell 'D49' '0.06657'
$ws.Range('E49').Value = '  -3.58%  '
Set-TextCell 'D50' '1.082'
$ws.Range('E50').Value = '  -4.58%  '
$ws.Range('E51').Value = '  +0.07%  '
